$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. data_lineup sheet: update roster stat cells (value tweaks)
# ---------------------------------------------------------------------------
$wsLineup = $wb.Worksheets.Item("data_lineup")

$wsLineup.Range("F2").Value = 35
$wsLineup.Range("H2").Value = 45
$wsLineup.Range("J2").Value = 34
$wsLineup.Range("M2").Value = 3

$wsLineup.Range("G3").Value = 69
$wsLineup.Range("J3").Value = 51

$wsLineup.Range("G4").Value = 96
$wsLineup.Range("J4").Value = 90

$wsLineup.Range("I5").Value = 54

$wsLineup.Range("F6").Value = 70
$wsLineup.Range("H6").Value = 77

$wsLineup.Range("M7").Value = 3

$wsLineup.Range("F8").Value = 99

$wsLineup.Range("F10").Value = 57
$wsLineup.Range("J10").Value = 88

$wsLineup.Range("E11").Value = 60
$wsLineup.Range("H11").Value = 76

$wsLineup.Range("J13").Value = 53
$wsLineup.Range("K13").Value = 46

# ---------------------------------------------------------------------------
# 2. Sheet selection / active tab state.
#    "as" was the active sheet before; now "data_lineup" is the active one,
#    scrolled/selected at cell U10. Update the no-longer-active "as" sheet's
#    remembered selection first, then activate data_lineup last so it ends
#    up as the truly selected/active tab.
# ---------------------------------------------------------------------------
$wsAs = $wb.Worksheets.Item("as")
$wsAs.Activate()
$wsAs.Range("D19").Select()

$wsLineup.Activate()
$wsLineup.Range("U10").Select()
